# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.188.48"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "3.073.28"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'521.71"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Value = "'135.40"
$ws.Range("E6").Value = "  -4.84%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.071.70"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("E9").Value = "  +4.46%  "
$ws.Range("D10").Value = "'7.29"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("D12").Value = "'0.401"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "3.602.24"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "57.235.90"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "3.069.02"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "'5.85"
$ws.Range("E19").Value = "  -4.37%  "
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").Value = "'349.90"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'68.97"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.165"
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("D28").Value = "0.0₃0863"
$ws.Range("E28").Value = "  -6.49%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "'5.82"
$ws.Range("E32").Value = "  -8.96%  "
$ws.Range("D33").Value = "'20.87"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").Value = "'4.83"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").Value = "'158.78"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  -5.34%  "
$ws.Range("D37").Value = "'5.99"
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("D38").Value = "'25.37"
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").Value = "'0.0655"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").Value = "'1.58"
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("D42").Value = "'4.05"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "2.402.63"
$ws.Range("D45").Value = "'36.58"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "3.112.33"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("D48").Value = "'0.0262"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").Value = "'0.943"
$ws.Range("E50").Value = "  -5.85%  "
$ws.Range("D51").Value = "'19.54"
$ws.Range("E51").Value = "  -5.32%  "
